# Update Gilgamesh_Profits market-price snapshot data (scheduled runner refresh).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1009.6667
$ws.Range("I28").Value = 328.5
$ws.Range("J28").Value = 2372
$ws.Range("K28").Value = 328.5
$ws.Range("L28").Value = 2372
$ws.Range("M28").Value = 156.5
$ws.Range("N28").Value = -3342
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
# Row 32
$ws.Range("H32").Value = 16669867
$ws.Range("I32").Value = 27779778
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 27779778
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -27779452
$ws.Range("N32").Value = -5652
# Row 62
$ws.Range("H62").Value = 3980.353
$ws.Range("I62").Value = 3497
$ws.Range("J62").Value = 4670.857
$ws.Range("K62").Value = 3497
$ws.Range("L62").Value = 4670.857
$ws.Range("M62").Value = -2873
$ws.Range("N62").Value = -5918.857
# Row 64
$ws.Range("H64").Value = 250003360
$ws.Range("J64").Value = 1000000000
$ws.Range("L64").Value = 1000000000
$ws.Range("N64").Value = -1000000496
# Row 65
$ws.Range("H65").Value = 3980.353
$ws.Range("I65").Value = 3497
$ws.Range("J65").Value = 4670.857
$ws.Range("K65").Value = 17485
$ws.Range("L65").Value = 23354.285
$ws.Range("M65").Value = -14365
$ws.Range("N65").Value = -29594.285
# Row 67
$ws.Range("H67").Value = 250003360
$ws.Range("J67").Value = 1000000000
$ws.Range("L67").Value = 1000000000
$ws.Range("N67").Value = -1000001716
# Row 112
$ws.Range("H112").Value = 1713.2222
$ws.Range("J112").Value = 1785.4706
$ws.Range("L112").Value = 5356.4118
$ws.Range("N112").Value = -7572.4118
# Row 113
$ws.Range("H113").Value = 3889.7778
$ws.Range("J113").Value = 4382.2
$ws.Range("L113").Value = 4382.2
$ws.Range("N113").Value = -10890.2

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 1394.0555
$ws.Range("I110").Value = 1269.0769
$ws.Range("K110").Value = 1269.0769
$ws.Range("M110").Value = 775.9231
# Row 132
$ws.Range("H132").Value = 3674.0476
$ws.Range("I132").Value = 3810.3333
$ws.Range("J132").Value = 3333.3333
$ws.Range("K132").Value = 11430.9999
$ws.Range("L132").Value = 9999.999899999999
$ws.Range("M132").Value = -8900.999899999999
$ws.Range("N132").Value = -15059.9999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 13001769
$ws.Range("I105").Value = 1112993.2
$ws.Range("J105").Value = 22728950
$ws.Range("K105").Value = 1112993.2
$ws.Range("L105").Value = 22728950
$ws.Range("M105").Value = -1111246.2
$ws.Range("N105").Value = -22732444
# Row 134
$ws.Range("H134").Value = 2092.7896
$ws.Range("I134").Value = 1443.3077
$ws.Range("K134").Value = 4329.9231
$ws.Range("M134").Value = -1794.9231

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2088.25
$ws.Range("I16").Value = 1958
$ws.Range("K16").Value = 1958
$ws.Range("M16").Value = -1671
# Row 22
$ws.Range("H22").Value = 300.25
$ws.Range("I22").Value = 317
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 317
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = 33
$ws.Range("N22").Value = -950
# Row 94
$ws.Range("H94").Value = 792.55554
$ws.Range("I94").Value = 596.75
$ws.Range("J94").Value = 949.2
$ws.Range("K94").Value = 596.75
$ws.Range("L94").Value = 949.2
$ws.Range("M94").Value = -145.75
$ws.Range("N94").Value = -1851.2
# Row 105
$ws.Range("H105").Value = 2157.7778
$ws.Range("I105").Value = 1774.2858
$ws.Range("K105").Value = 1774.2858
$ws.Range("M105").Value = -27.28580000000011
# Row 113
$ws.Range("H113").Value = 2088.25
$ws.Range("I113").Value = 1958
$ws.Range("K113").Value = 1958
$ws.Range("M113").Value = 212
# Row 122
$ws.Range("H122").Value = 408.9524
$ws.Range("I122").Value = 374.41177
$ws.Range("K122").Value = 1123.23531
$ws.Range("M122").Value = 1326.76469
# Row 134
$ws.Range("H134").Value = 4423.16
$ws.Range("I134").Value = 4999.3335
$ws.Range("K134").Value = 14998.0005
$ws.Range("M134").Value = -12463.0005

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 183.16667
$ws.Range("J12").Value = 272
$ws.Range("L12").Value = 816
$ws.Range("N12").Value = -1162
# Row 98
$ws.Range("H98").Value = 890
$ws.Range("J98").Value = 900
$ws.Range("L98").Value = 2700
$ws.Range("N98").Value = -5696
# Row 131
$ws.Range("H131").Value = 6324.769
$ws.Range("I131").Value = 12491.8
$ws.Range("J131").Value = 2470.375
$ws.Range("K131").Value = 37475.39999999999
$ws.Range("L131").Value = 7411.125
$ws.Range("M131").Value = -32435.39999999999
$ws.Range("N131").Value = -17491.125
# Row 141
$ws.Range("H141").Value = 8055.6875
$ws.Range("I141").Value = 2299.3076
$ws.Range("K141").Value = 6897.9228
$ws.Range("M141").Value = -1717.9228

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 51
$ws.Range("H51").Value = 49998.5
$ws.Range("J51").Value = 49998.5
$ws.Range("L51").Value = 49998.5
$ws.Range("N51").Value = -51016.5
# Row 97
$ws.Range("H97").Value = 1808.9
$ws.Range("I97").Value = 1769.3529
$ws.Range("K97").Value = 1769.3529
$ws.Range("M97").Value = -1273.3529
# Row 102
$ws.Range("H102").Value = 21319.9
$ws.Range("I102").Value = 1531.3334
$ws.Range("J102").Value = 51002.75
$ws.Range("K102").Value = 1531.3334
$ws.Range("L102").Value = 51002.75
$ws.Range("M102").Value = 90.66660000000002
$ws.Range("N102").Value = -54246.75
# Row 113
$ws.Range("H113").Value = 2082.8572
$ws.Range("I113").Value = 1775.75
$ws.Range("K113").Value = 1775.75
$ws.Range("M113").Value = 394.25
# Row 122
$ws.Range("H122").Value = 2751212.5
$ws.Range("I122").Value = 5131449
$ws.Range("J122").Value = 4785.769
$ws.Range("K122").Value = 15394347
$ws.Range("L122").Value = 14357.307
$ws.Range("M122").Value = -15391897
$ws.Range("N122").Value = -19257.307

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4256.5
$ws.Range("I7").Value = 3876
$ws.Range("K7").Value = 3876
$ws.Range("M7").Value = -3764
# Row 61
$ws.Range("H61").Value = 1532.3334
$ws.Range("I61").Value = 1532.3334
$ws.Range("K61").Value = 1532.3334
$ws.Range("M61").Value = -1330.3334
# Row 113
$ws.Range("H113").Value = 1532.3334
$ws.Range("I113").Value = 1532.3334
$ws.Range("K113").Value = 1532.3334
$ws.Range("M113").Value = 637.6666
# Row 122
$ws.Range("H122").Value = 16913.834
$ws.Range("I122").Value = 11249.5
$ws.Range("K122").Value = 33748.5
$ws.Range("M122").Value = -31298.5
# Row 126
$ws.Range("H126").Value = 4256.5
$ws.Range("I126").Value = 3876
$ws.Range("K126").Value = 11628
$ws.Range("M126").Value = -9158
# Row 132
$ws.Range("H132").Value = 4997.467
$ws.Range("I132").Value = 3877.4
$ws.Range("K132").Value = 11632.2
$ws.Range("M132").Value = -9102.200000000001
# Row 136
$ws.Range("H136").Value = 1763.8049
$ws.Range("I136").Value = 1271.4736
$ws.Range("K136").Value = 3814.4208
$ws.Range("M136").Value = -1264.4208

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 43
$ws.Range("H43").Value = 22788.166
$ws.Range("I43").Value = 22499.5
$ws.Range("J43").Value = 22932.5
$ws.Range("K43").Value = 22499.5
$ws.Range("L43").Value = 22932.5
$ws.Range("M43").Value = -22350.5
$ws.Range("N43").Value = -23230.5
# Row 96
$ws.Range("H96").Value = 1241
$ws.Range("I96").Value = 1241
$ws.Range("K96").Value = 1241
$ws.Range("M96").Value = 132
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
# Row 113
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830
# Row 126
$ws.Range("H126").Value = 1300
$ws.Range("J126").Value = 1500
$ws.Range("L126").Value = 4500
$ws.Range("N126").Value = -9440
# Row 132
$ws.Range("H132").Value = 1365.3611
$ws.Range("I132").Value = 1289.0358
$ws.Range("K132").Value = 3867.1074
$ws.Range("M132").Value = -1337.1074
# Row 136
$ws.Range("H136").Value = 4049.1614
$ws.Range("I136").Value = 2086.862
$ws.Range("K136").Value = 6260.586
$ws.Range("M136").Value = -3710.586

